$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74; existing rows 74-123 shift down to 75-124.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new data record.
$ws.Cells.Item(74, 1).Value = 7
$ws.Cells.Item(74, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(74, 3).Value = "Ñuble"
$ws.Cells.Item(74, 4).Value = 44582
$ws.Cells.Item(74, 5).Value = 16
$ws.Cells.Item(74, 6).Value = 100112045
$ws.Cells.Item(74, 7).Value = "Zapallo"
$ws.Cells.Item(74, 8).Value = "Camote"
$ws.Cells.Item(74, 9).Value = "1a nueva(o)"
$ws.Cells.Item(74, 10).Value = 300
$ws.Cells.Item(74, 11).Value = 300
$ws.Cells.Item(74, 12).Value = 350
$ws.Cells.Item(74, 13).Value = 325
$ws.Cells.Item(74, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(74, 15).Value = "Región del Maule"
$ws.Cells.Item(74, 16).Value = 325
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
